$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from H1 (existing header cell) to new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data rows
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 3

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 5

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 5
